$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.970.23"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "1.762.60"
$ws.Range("E3").Value = "  -3.48%  "
$ws.Range("E4").Value = "  +0.79%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "339.37"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3771"
$ws.Range("D7").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3375"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.14"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.132"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -6.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07217"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -5.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.55"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.207"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.205"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("D16").Value = "1.762.34"
$ws.Range("E16").Value = "  -3.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001054"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -5.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06572"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.59"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -5.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9998"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.96"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.273"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.98%  "
$ws.Range("D23").Value = "27.979.50"
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.70"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -9.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.396"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.15"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.79"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -8.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.330"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -10.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.285"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -16.40%  "
$ws.Range("D30").Value = "1.964.64"
$ws.Range("E30").Value = "  -3.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "131.31"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.016"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.833"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -6.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08792"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.24"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -8.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02347"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6595"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06220"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.148"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -7.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2113"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -5.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.211"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.461"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -9.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.062"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9992"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.73"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -6.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.838"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6051"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -7.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.01"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.012"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -7.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07229"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.183"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.92%  "
